$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 3, 4, 5 data gets cyclically rotated:
#   new row 3 <- old row 5
#   new row 4 <- old row 3
#   new row 5 <- old row 4
# Only columns D, M, N, O, P, S actually change values (others stay identical
# between the rotated rows), so we just set the final target values directly.

# Row 3 (becomes old row 5's values)
$ws.Range("D3").Value = 44981
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 25000
$ws.Range("S3").Value = 3125

# Row 4 (becomes old row 3's values)
$ws.Range("D4").Value = 44973
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 28000
$ws.Range("S4").Value = 3500

# Row 5 (becomes old row 4's values)
$ws.Range("D5").Value = 44971
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 28000
$ws.Range("S5").Value = 3500
